# "Activity Tracker" workbook update:
# add a new row (#11) to the tracker table noting the work is finished.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New record: Sr.No 11, Activities = "Its finished now"
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Its finished now"

# Match the wrap-text style used by the other "Activities" column cells (B2:B11)
$ws.Cells.Item(12, 2).WrapText = $true

# Leave the selection on the newly entered cell, as the author did
$ws.Range("B12").Select()
